# Fruta / hortaliza, semanal
# Insert a new weekly record at row 48, pushing the existing data rows
# (old rows 48-100) down by one row (new rows 49-101).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 48.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new record.
$ws.Cells.Item(48, 1).Value = 10
$ws.Cells.Item(48, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(48, 3).Value = "La Araucanía"
$ws.Cells.Item(48, 4).Value = 44539
$ws.Cells.Item(48, 5).Value = 9
$ws.Cells.Item(48, 6).Value = 100112031
$ws.Cells.Item(48, 7).Value = "Poroto verde"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 150
$ws.Cells.Item(48, 11).Value = 23000
$ws.Cells.Item(48, 12).Value = 23000
$ws.Cells.Item(48, 13).Value = 23000
$ws.Cells.Item(48, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(48, 15).Value = "Región del Maule"
$ws.Cells.Item(48, 16).Value = 920
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
